$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- Text content updates ----
$ws.Range("B3").Value = "I0"
$ws.Range("B4").Value = "I1"
$ws.Range("B5").Value = "I2"
$ws.Range("B6").Value = "I3"

$ws.Range("F7").Value = "I0, model::execute(I0), check=0 ->  activate check"
$ws.Range("F8").Value = " check state after I0, model::get_reg , I1, model::execute(I1)"
$ws.Range("F9").Value = " check state after I1, I2, model::execute(I2)"
$ws.Range("F10").Value = " check state after I2, I3, model::execute(I3)"

$ws.Range("B7").ClearContents()
$ws.Range("B8").ClearContents()
$ws.Range("B9").ClearContents()
$ws.Range("B10").ClearContents()

# ---- Style updates ----
# B7: gray fill + border, alignment flag touched (matches new cellXfs index 7)
$ws.Range("C8").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B7").WrapText = $false

# B8/B9/B10: gray fill + border + wrap text (matches new cellXfs index 8)
$ws.Range("C8").Copy()
$ws.Range("B8").PasteSpecial(-4122)
$ws.Range("B8").WrapText = $true

$ws.Range("C8").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").WrapText = $true

$ws.Range("C8").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B10").WrapText = $true

# F8/F9/F10: border + wrap text (matches existing cellXfs index 4)
$ws.Range("F7").Copy()
$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F7").Copy()
$ws.Range("F10").PasteSpecial(-4122)

# B11: gray fill + border (matches existing cellXfs index 5)
$ws.Range("C8").Copy()
$ws.Range("B11").PasteSpecial(-4122)

# B12: new row/cell, gray fill only, no border (matches new cellXfs index 9)
$ws.Range("C8").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Borders.LineStyle = -4142

$excel.CutCopyMode = $false

# ---- Row heights ----
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30

# ---- Column width ----
$ws.Columns.Item(2).ColumnWidth = 6.14

# ---- Selection ----
$ws.Range("J7").Select()
